# Bugfixed the naive forecaster component module
#
# The naive forecaster's YoY forecast vector sheet had an off-by-one row
# (the first data row was a stray/incomplete record) and the C/E forecast
# columns were computed with a bug. This script:
#   1. Deletes the obsolete first data row (old row 2), shifting all other
#      data rows up by one (dates/years in A/B/D realign correctly).
#   2. Re-writes the recomputed forecast columns C (y_0_forecast) and
#      E (y_1_forecast) with the corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the stale leading data row; everything below shifts up one row.
$ws.Rows(2).Delete()

# 2) Corrected forecast values (post-bugfix) for rows 2..52.
$cValues = @(
    $null, $null, -0.9921462019007898, -1.120856461469888, 0.2467309912830284,
    0.2946875655135184, 1.160201558804674, 1.242549344471122, 1.048604932640185,
    1.180122656701199, 1.06837811337479, 1.381744454873757, 1.435208340819005,
    1.404039427736437, 1.577589817310243, 1.593617458167307, 1.979569114089963,
    1.885212754467758, 1.815212363528707, 1.877372574190161, 2.140635848901895,
    2.284026378382942, 2.336516087993035, 2.279995067217899, 1.404530461900833,
    1.264761787657309, 1.17909021197069, 1.025257057800411, 0.5345697479163913,
    0.824608016336259, -1.788000783651811, -1.788000783651811, -2.680286313062752,
    -1.4191429117966, -1.098964423305859, -1.098964423305859, 1.514644056931957,
    1.896944139870205, 1.916393754370604, 1.916393754370604, -0.8557279162653919,
    -0.7016063587211741, -0.7359525160776204, -0.7359525160776204, -0.1316183744203947,
    -0.1754728623905355, -0.187152549496028, -0.187152549496028, 0.3903331526556864,
    0.5695821893874298, 0.6150340712028246
)

$eValues = @(
    1.782259294303912, -0.08289353495386509, -0.3230348957779294, -0.363786394693788,
    -0.06959526544320083, -0.2139598932957232, 0.2932139896134167, 0.903223459378788,
    1.078804187516891, 1.31837503023402, 1.25598608434605, 1.658305347589661,
    1.407107513712802, 1.488472133572305, 1.464859320654099, 1.644157643645183,
    1.639881111696151, 1.586470485311331, 1.806931013599544, 1.974604558490256,
    1.99288634244883, 2.119133965447961, 2.162438527487853, 1.93172124148786,
    1.887821778955101, 1.798687504247187, 1.636329093826605, 1.008270799755984,
    1.247274949485733, 1.395219579261608, -0.02261741485058977, -1.119700950349478,
    -2.013357217277445, -0.1125839228000469, 1.055324027461602, 0.5759895884974942,
    -0.06175132635745095, 0.5116467003986713, 0.4136280550221194, 0.3530477102890783,
    -0.2004689067778398, 0.2942159770784825, 0.6923809915882817, -0.01286797263981843,
    -0.09133135081734745, -0.0331361487157622, -0.3126391654689975, -0.1152140120150968,
    -0.2264357368625403, 0.316149716722669, 0.6473947787101642
)

for ($i = 0; $i -lt 51; $i++) {
    $row = $i + 2
    $cVal = $cValues[$i]
    if ($null -eq $cVal) {
        $ws.Cells.Item($row, 3).ClearContents()
    } else {
        $ws.Cells.Item($row, 3).Value = $cVal
    }
    $ws.Cells.Item($row, 5).Value = $eValues[$i]
}
